$wb = $excel.ActiveWorkbook

# Sheet 1: 台指期換倉成本計算 - insert new day's row at row 2
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()
$ws1.Cells.Item(2,1).Value = "日期：2021/12/01"
$ws1.Cells.Item(2,2).Value = "202201"
$ws1.Cells.Item(2,3).Value = 17536
$ws1.Cells.Item(2,4).Value = 7479
$ws1.Cells.Item(2,5).Value = 8592640
$ws1.Cells.Item(2,6).Value = 17630

# Sheet 2: 散戶多空力道 - insert new day's row at row 2
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Insert()
$ws2.Cells.Item(2,1).Value = "日期：2021/12/01"
$ws2.Cells.Item(2,2).Value = 0.01

# Sheet 3: 三大法人買賣金額 - insert new day's row at row 2
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(2).Insert()
$ws3.Cells.Item(2,1).Value = "110年12月01日"
$ws3.Cells.Item(2,2).Value = 84.34
$ws3.Cells.Item(2,3).Value = 73.89

# Sheet 4: 大盤多空點位 - insert new day's row at row 2
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Insert()
$ws4.Cells.Item(2,1).Value = "110年12月01日"
$ws4.Cells.Item(2,2).Value = 17487.19

# Sheet 5: 期貨大額交易人未沖銷部位 - insert new day's row at row 2
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows.Item(2).Insert()
$ws5.Cells.Item(2,1).Value = "2021/12/01"
$ws5.Cells.Item(2,2).Value = 49411
$ws5.Cells.Item(2,3).Value = 54693
$ws5.Cells.Item(2,4).Value = 288
$ws5.Cells.Item(2,5).Value = -648
$ws5.Cells.Item(2,6).Value = 27067
$ws5.Cells.Item(2,7).Value = 48641
$ws5.Cells.Item(2,8).Value = 1251
$ws5.Cells.Item(2,9).Value = 1031
$ws5.Cells.Item(2,10).Value = -21574
$ws5.Cells.Item(2,11).Value = 220
$ws5.Cells.Item(2,12).Value = -963
$ws5.Cells.Item(2,13).Value = -1679
$ws5.Cells.Item(2,14).Value = 716
